$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033095748054874
$ws.Cells.Item(2, 4).Value = 1.035837824644373
$ws.Cells.Item(2, 5).Value = 1.032367888000352
$ws.Cells.Item(2, 6).Value = 1.041264067148589
$ws.Cells.Item(2, 9).Value = 1.033257053570356
$ws.Cells.Item(2, 10).Value = 1.038222354002994
$ws.Cells.Item(2, 11).Value = 1.038633284321468
$ws.Cells.Item(2, 12).Value = 1.035173321254598
$ws.Cells.Item(2, 13).Value = 1.044044074934244
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034725863497439
$ws.Cells.Item(3, 4).Value = 1.037074053980287
$ws.Cells.Item(3, 5).Value = 1.033776928272614
$ws.Cells.Item(3, 6).Value = 1.042974392707843
$ws.Cells.Item(3, 9).Value = 1.033644668577636
$ws.Cells.Item(3, 10).Value = 1.039491589352489
$ws.Cells.Item(3, 11).Value = 1.039677857957511
$ws.Cells.Item(3, 12).Value = 1.036389513817423
$ws.Cells.Item(3, 13).Value = 1.045562631427253
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.035778243914784
$ws.Cells.Item(4, 4).Value = 1.037871668507612
$ws.Cells.Item(4, 5).Value = 1.034686703648659
$ws.Cells.Item(4, 6).Value = 1.044078879040007
$ws.Cells.Item(4, 9).Value = 1.033893043715278
$ws.Cells.Item(4, 10).Value = 1.040310163712612
$ws.Cells.Item(4, 11).Value = 1.040350909144308
$ws.Cells.Item(4, 12).Value = 1.037173994081588
$ws.Cells.Item(4, 13).Value = 1.046542583310783
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.036220098941274
$ws.Cells.Item(5, 4).Value = 1.038206441095392
$ws.Cells.Item(5, 5).Value = 1.035068711459252
$ws.Cells.Item(5, 6).Value = 1.044542688855022
$ws.Cells.Item(5, 9).Value = 1.03399688034239
$ws.Cells.Item(5, 10).Value = 1.04065365384867
$ws.Cells.Item(5, 11).Value = 1.040633183284157
$ws.Cells.Item(5, 12).Value = 1.037503205066693
$ws.Cells.Item(5, 13).Value = 1.046953930698858
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036294255519628
$ws.Cells.Item(6, 4).Value = 1.038262619186853
$ws.Cells.Item(6, 5).Value = 1.035132825444455
$ws.Cells.Item(6, 6).Value = 1.044620534653301
$ws.Cells.Item(6, 9).Value = 1.034014281036542
$ws.Cells.Item(6, 10).Value = 1.040711290177688
$ws.Cells.Item(6, 11).Value = 1.040680538898501
$ws.Cells.Item(6, 12).Value = 1.03755844701284
$ws.Cells.Item(6, 13).Value = 1.047022961450647
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.035784150209493
$ws.Cells.Item(7, 4).Value = 1.037876143883803
$ws.Cells.Item(7, 5).Value = 1.034691809857145
$ws.Cells.Item(7, 6).Value = 1.044085078500672
$ws.Cells.Item(7, 9).Value = 1.033894433460862
$ws.Cells.Item(7, 10).Value = 1.040314755940798
$ws.Cells.Item(7, 11).Value = 1.040354683553151
$ws.Cells.Item(7, 12).Value = 1.037178395299682
$ws.Cells.Item(7, 13).Value = 1.046548082192398
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033647160161469
$ws.Cells.Item(8, 4).Value = 1.036256097028968
$ws.Cells.Item(8, 5).Value = 1.032844492881393
$ws.Cells.Item(8, 6).Value = 1.041842544563812
$ws.Cells.Item(8, 9).Value = 1.033388556616968
$ws.Cells.Item(8, 10).Value = 1.038651863868986
$ws.Cells.Item(8, 11).Value = 1.038986899184738
$ws.Cells.Item(8, 12).Value = 1.035584857458016
$ws.Cells.Item(8, 13).Value = 1.044557834980559
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029862485496266
$ws.Cells.Item(9, 4).Value = 1.033383321459145
$ws.Cells.Item(9, 5).Value = 1.029573790681578
$ws.Cells.Item(9, 6).Value = 1.037873433578813
$ws.Cells.Item(9, 9).Value = 1.032478317804452
$ws.Cells.Item(9, 10).Value = 1.03570049581953
$ws.Cells.Item(9, 11).Value = 1.03655445887289
$ws.Cells.Item(9, 12).Value = 1.032757475077431
$ws.Cells.Item(9, 13).Value = 1.041029912062324
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027325815932184
$ws.Cells.Item(10, 4).Value = 1.031455489934723
$ws.Cells.Item(10, 5).Value = 1.027382312206214
$ws.Cells.Item(10, 6).Value = 1.035214815411137
$ws.Cells.Item(10, 9).Value = 1.031858628214961
$ws.Cells.Item(10, 10).Value = 1.033718117961167
$ws.Cells.Item(10, 11).Value = 1.034917404475412
$ws.Cells.Item(10, 12).Value = 1.030858985989487
$ws.Cells.Item(10, 13).Value = 1.038663221705196
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026224024956076
$ws.Cells.Item(11, 4).Value = 1.030617602143003
$ws.Cells.Item(11, 5).Value = 1.026430633345172
$ws.Cells.Item(11, 6).Value = 1.034060451929173
$ws.Cells.Item(11, 9).Value = 1.031587197210531
$ws.Cells.Item(11, 10).Value = 1.032856084945111
$ws.Cells.Item(11, 11).Value = 1.034204774506837
$ws.Cells.Item(11, 12).Value = 1.03003357847148
$ws.Cells.Item(11, 13).Value = 1.037634760636036
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025814245696935
$ws.Cells.Item(12, 4).Value = 1.030305893898078
$ws.Cells.Item(12, 5).Value = 1.026076712100224
$ws.Cells.Item(12, 6).Value = 1.033631179185867
$ws.Cells.Item(12, 9).Value = 1.03148590554813
$ws.Cells.Item(12, 10).Value = 1.032535327662817
$ws.Cells.Item(12, 11).Value = 1.033939495491647
$ws.Cells.Item(12, 12).Value = 1.02972647151333
$ws.Cells.Item(12, 13).Value = 1.037252179239908
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025902168736774
$ws.Cells.Item(13, 4).Value = 1.030372778231901
$ws.Cells.Item(13, 5).Value = 1.026152648843107
$ws.Cells.Item(13, 6).Value = 1.033723282102078
$ws.Cells.Item(13, 9).Value = 1.031507654294006
$ws.Cells.Item(13, 10).Value = 1.03260415670719
$ws.Cells.Item(13, 11).Value = 1.033996424981818
$ws.Cells.Item(13, 12).Value = 1.029792370423936
$ws.Cells.Item(13, 13).Value = 1.037334270049359
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.02619016327123
$ws.Cells.Item(14, 4).Value = 1.030591846080718
$ws.Cells.Item(14, 5).Value = 1.026401386848785
$ws.Cells.Item(14, 6).Value = 1.034024978210795
$ws.Cells.Item(14, 9).Value = 1.03157883402431
$ws.Cells.Item(14, 10).Value = 1.032829582535119
$ws.Cells.Item(14, 11).Value = 1.0341828582763
$ws.Cells.Item(14, 12).Value = 1.03000820345915
$ws.Cells.Item(14, 13).Value = 1.0376031479492
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026367536192252
$ws.Cells.Item(15, 4).Value = 1.03072675718279
$ws.Cells.Item(15, 5).Value = 1.026554585833761
$ws.Cells.Item(15, 6).Value = 1.034210797649698
$ws.Cells.Item(15, 9).Value = 1.031622627786949
$ws.Cells.Item(15, 10).Value = 1.032968400300011
$ws.Cells.Item(15, 11).Value = 1.034297649304603
$ws.Cells.Item(15, 12).Value = 1.030141116903538
$ws.Cells.Item(15, 13).Value = 1.037768737186225
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027398865356491
$ws.Cells.Item(16, 4).Value = 1.031511031112012
$ws.Cells.Item(16, 5).Value = 1.027445412986481
$ws.Cells.Item(16, 6).Value = 1.035291358739342
$ws.Cells.Item(16, 9).Value = 1.031876576501219
$ws.Cells.Item(16, 10).Value = 1.033775250373486
$ws.Cells.Item(16, 11).Value = 1.034964619050841
$ws.Cells.Item(16, 12).Value = 1.030913694095551
$ws.Cells.Item(16, 13).Value = 1.038731398795248
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028044871220771
$ws.Cells.Item(17, 4).Value = 1.032002141969011
$ws.Cells.Item(17, 5).Value = 1.02800345982421
$ws.Cells.Item(17, 6).Value = 1.035968308209439
$ws.Cells.Item(17, 9).Value = 1.032035038584023
$ws.Cells.Item(17, 10).Value = 1.034280380790188
$ws.Cells.Item(17, 11).Value = 1.035381974361457
$ws.Cells.Item(17, 12).Value = 1.031397407177848
$ws.Cells.Item(17, 13).Value = 1.03933425897727
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028421348822637
$ws.Cells.Item(18, 4).Value = 1.03228829788622
$ws.Cells.Item(18, 5).Value = 1.028328694147911
$ws.Cells.Item(18, 6).Value = 1.036362857121223
$ws.Cells.Item(18, 9).Value = 1.032127167887168
$ws.Cells.Item(18, 10).Value = 1.034574663390821
$ws.Cells.Item(18, 11).Value = 1.035625047099875
$ws.Cells.Item(18, 12).Value = 1.031679226411375
$ws.Cells.Item(18, 13).Value = 1.039685544116076
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028549662988278
$ws.Cells.Item(19, 4).Value = 1.032385818951487
$ws.Cells.Item(19, 5).Value = 1.028439546053953
$ws.Cells.Item(19, 6).Value = 1.036497336958106
$ws.Cells.Item(19, 9).Value = 1.032158531016939
$ws.Cells.Item(19, 10).Value = 1.034674946888344
$ws.Cells.Item(19, 11).Value = 1.03570786724053
$ws.Cells.Item(19, 12).Value = 1.03177526511955
$ws.Cells.Item(19, 13).Value = 1.039805263815975
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027975594803012
$ws.Cells.Item(20, 4).Value = 1.031949481601222
$ws.Cells.Item(20, 5).Value = 1.027943614168881
$ws.Cells.Item(20, 6).Value = 1.035895709519934
$ws.Cells.Item(20, 9).Value = 1.032018068053972
$ws.Cells.Item(20, 10).Value = 1.034226221509052
$ws.Cells.Item(20, 11).Value = 1.035337233758533
$ws.Cells.Item(20, 12).Value = 1.03134554276884
$ws.Cells.Item(20, 13).Value = 1.03926961439939
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026105370671823
$ws.Cells.Item(21, 4).Value = 1.03052734937486
$ws.Cells.Item(21, 5).Value = 1.026328151540666
$ws.Cells.Item(21, 6).Value = 1.033936149897493
$ws.Cells.Item(21, 9).Value = 1.031557886378629
$ws.Cells.Item(21, 10).Value = 1.032763215793968
$ws.Cells.Item(21, 11).Value = 1.034127974308627
$ws.Cells.Item(21, 12).Value = 1.029944660280966
$ws.Cells.Item(21, 13).Value = 1.037523985823793
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024926441419925
$ws.Cells.Item(22, 4).Value = 1.029630419603167
$ws.Cells.Item(22, 5).Value = 1.025309979356702
$ws.Cells.Item(22, 6).Value = 1.032701249950322
$ws.Cells.Item(22, 9).Value = 1.031265829933341
$ws.Cells.Item(22, 10).Value = 1.031840121238568
$ws.Cells.Item(22, 11).Value = 1.033364324616097
$ws.Cells.Item(22, 12).Value = 1.029060891485159
$ws.Cells.Item(22, 13).Value = 1.036423162950954
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025551707736269
$ws.Cells.Item(23, 4).Value = 1.030106165831835
$ws.Cells.Item(23, 5).Value = 1.025849969400086
$ws.Cells.Item(23, 6).Value = 1.033356168810352
$ws.Cells.Item(23, 9).Value = 1.031420913993506
$ws.Cells.Item(23, 10).Value = 1.032329782540734
$ws.Cells.Item(23, 11).Value = 1.033769469559322
$ws.Cells.Item(23, 12).Value = 1.029529679953493
$ws.Cells.Item(23, 13).Value = 1.037007045311349
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028006898843596
$ws.Cells.Item(24, 4).Value = 1.03197327748399
$ws.Cells.Item(24, 5).Value = 1.027970656663211
$ws.Cells.Item(24, 6).Value = 1.035928514682806
$ws.Cells.Item(24, 9).Value = 1.03202573722976
$ws.Cells.Item(24, 10).Value = 1.0342506948428
$ws.Cells.Item(24, 11).Value = 1.035357451234523
$ws.Cells.Item(24, 12).Value = 1.031368979059058
$ws.Cells.Item(24, 13).Value = 1.03929882559223
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030843247153842
$ws.Cells.Item(25, 4).Value = 1.034128193999135
$ws.Cells.Item(25, 5).Value = 1.03042124238626
$ws.Cells.Item(25, 6).Value = 1.03890169638392
$ws.Cells.Item(25, 9).Value = 1.032715888075868
$ws.Cells.Item(25, 10).Value = 1.036466060332154
$ws.Cells.Item(25, 11).Value = 1.037185986579461
$ws.Cells.Item(25, 12).Value = 1.033490772079454
$ws.Cells.Item(25, 13).Value = 1.041944509716573

Write-Host "Updated vm_pu values for Case_5_72 (380 kV) bus voltage results."
